$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds numeric-looking text that must stay text (as in
# the source XML: t="inlineStr"). A plain Range.Value assignment of a
# numeric-looking string gets auto-coerced to a real number by Excel, so
# instead we write a text formula that evaluates to the exact string, then
# convert it in place to a literal value (Copy + PasteSpecial values-only).
# That keeps the cell's type as text without leaving a formula behind and
# without permanently changing the cell's NumberFormat/style.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.Formula = '="' + $val + '"'
    $rng.Copy()
    $rng.PasteSpecial(-4163)
}

# --- Simple Price (column D) updates ---
Set-TextValue "D2"  "264.54"
Set-TextValue "D3"  "22.71"
Set-TextValue "D4"  "6.220"
Set-TextValue "D5"  "0.06125"
Set-TextValue "D6"  "3.542"
Set-TextValue "D7"  "6.736"
Set-TextValue "D8"  "1.386"
Set-TextValue "D9"  "0.8157"
Set-TextValue "D10" "0.1600"
Set-TextValue "D11" "0.08224"
Set-TextValue "D13" "0.03163"
Set-TextValue "D14" "0.09265"
Set-TextValue "D15" "3.895"
Set-TextValue "D16" "0.001696"
Set-TextValue "D17" "0.04846"
Set-TextValue "D18" "0.0006255"
Set-TextValue "D19" "0.006204"

# --- Rows 20/21: coin order swapped (Coin, Link, Price, Volume) ---
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D20" "0.001101"
$ws.Range("E20").Value = "19BitKanKAN"

$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D21" "0.003201"
$ws.Range("E21").Value = "20HotbitTokenHTB"

# --- More Price (column D) updates ---
Set-TextValue "D23" "3.698"
Set-TextValue "D24" "2.254"
Set-TextValue "D25" "0.3406"
Set-TextValue "D26" "0.1271"
Set-TextValue "D27" "0.0002688"
Set-TextValue "D40" "0.04650"

# --- Rows 41/42/43: coin order rotated (Coin, Link, Price, Volume) ---
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.007213"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1127"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003136"
$ws.Range("E43").Value = "42CEJICEJI"

# --- More Price (column D) updates ---
Set-TextValue "D44" "0.01035"
Set-TextValue "D45" "0.00006162"
Set-TextValue "D46" "0.00000000752"
Set-TextValue "D47" "0.7518"
Set-TextValue "D48" "0.1665"
Set-TextValue "D49" "0.00002105"
Set-TextValue "D50" "0.01243"

$excel.CutCopyMode = $false
